$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4 record - RollNumber first so shared string ordering matches source
$ws.Range("B4").Value = "HE130604"

# Add Fullname for row 3 (previously missing)
$ws.Range("C3").Value = "Pham Thanh Ha"

$ws.Range("A4").Value = 3
$ws.Range("E4").Value = "Hà Nội"
$ws.Range("F4").Value = "Female"
$ws.Range("G4").Value = "Kinh"
$ws.Range("H4").Value = "Viet Nam"
$ws.Range("J3").Copy()
$ws.Range("J4").PasteSpecial(-4122)
$ws.Range("J4").Value = 54850
$ws.Range("K4").Value = "Bad"
$ws.Range("L4").Value = "QD@@"
$ws.Range("M4").Value = "$$"
$ws.Range("N4").Value = "ĐH 23232"
$ws.Range("I4").Value = "ES2"

$ws.Range("I6").Select() | Out-Null
